$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E1").Value = "Tiempo"
$ws.Range("E2").Value = 8
$ws.Range("E3").Value = 16
$ws.Range("E4").Value = 400

$ws.Range("C17").Select()
